$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# 1. Update the "Ready for handoff" status text (for the 9c2dbcd1... row, row 8)
#    to "Handback transform failed" everywhere it is shown:
#    Overview!E8, Overview!F8, zh-cn!C8, de-de!C8
$newStatus = "Handback transform failed"
$wsOverview.Range("E8").Value = $newStatus
$wsOverview.Range("F8").Value = $newStatus
$wsZhCn.Range("C8").Value = $newStatus
$wsDeDe.Range("C8").Value = $newStatus

# 2. Populate the Error Detail column (R8) on the zh-cn and de-de sheets
#    with the handback/handoff file name mismatch message.
$wsZhCn.Range("R8").Value = "Handback file name: rw0ynko3.3en is different with handoff file name: 9c2dbcd1-f134-4efa-a5bd-9cf87984f5e6.5495184ae7eb522b115bf9a0370f15077f2133a9.zh-cn."
$wsDeDe.Range("R8").Value = "Handback file name: rw0ynko3.3en is different with handoff file name: 9c2dbcd1-f134-4efa-a5bd-9cf87984f5e6.5495184ae7eb522b115bf9a0370f15077f2133a9.de-de."
